{"js": "// Update the date line and the 25 multiplication answers in the practice\n// table. Each edit rewrites only the text content of the run/cell it\n// targets, so existing run formatting (fonts, size, alignment, etc.) is\n// left untouched.\n\n// 1) Date heading paragraph: \"2025-10-07 Tuesday\" -> \"2025-10-08 Wednesday\"\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items.find(p => p.text.trim() === \"2025-10-07 Tuesday\");\nif (dateParagraph) {\n  dateParagraph.getRange().insertText(\"2025-10-08 Wednesday\", \"Replace\");\n}\n\n// 2) The 25 answer cells in the practice table, in row-major (reading) order.\nconst newAnswers = [\n  \"451\u00d73=1353\", \"438\u00d79=3942\", \"596\u00d74=2384\", \"222\u00d73=666\", \"680\u00d74=2720\",\n  \"956\u00d75=4780\", \"114\u00d78=912\", \"918\u00d78=7344\", \"107\u00d77=749\", \"520\u00d73=1560\",\n  \"795\u00d77=5565\", \"192\u00d76=1152\", \"850\u00d74=3400\", \"200\u00d76=1200\", \"441\u00d75=2205\",\n  \"732\u00d79=6588\", \"444\u00d74=1776\", \"542\u00d76=3252\", \"773\u00d75=3865\", \"149\u00d73=447\",\n  \"167\u00d79=1503\", \"171\u00d72=342\", \"518\u00d75=2590\", \"902\u00d72=1804\", \"157\u00d76=942\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nlet answerIndex = 0;\nfor (let row = 0; row < table.rowCount && answerIndex < newAnswers.length; row++) {\n  for (let col = 0; col < 5 && answerIndex < newAnswers.length; col++) {\n    const cell = table.getCellOrNullObject(row, col);\n    cell.load(\"value\");\n    await context.sync();\n    if (cell.isNullObject) continue;\n    if (cell.value && cell.value.trim() !== \"\") {\n      cell.value = newAnswers[answerIndex];\n      answerIndex++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 multiplication answers in the practice\n# table. Each call finds the exact old expression and swaps in the new one,\n# so existing run formatting (fonts, size, alignment) is left untouched.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText([string]$findText, [string]$replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # wdFindContinue=1, wdReplaceOne=1 -> replace only the single exact match\n    $find.Execute(\n        $findText,   # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $replaceText,# ReplaceWith\n        1            # Replace (wdReplaceOne)\n    ) | Out-Null\n}\n\nReplace-ExactText \"2025-10-07 Tuesday\" \"2025-10-08 Wednesday\"\nReplace-ExactText \"481\u00d79=4329\" \"451\u00d73=1353\"\nReplace-ExactText \"436\u00d72=872\" \"438\u00d79=3942\"\nReplace-ExactText \"665\u00d76=3990\" \"596\u00d74=2384\"\nReplace-ExactText \"315\u00d77=2205\" \"222\u00d73=666\"\nReplace-ExactText \"296\u00d75=1480\" \"680\u00d74=2720\"\nReplace-ExactText \"858\u00d73=2574\" \"956\u00d75=4780\"\nReplace-ExactText \"940\u00d73=2820\" \"114\u00d78=912\"\nReplace-ExactText \"481\u00d78=3848\" \"918\u00d78=7344\"\nReplace-ExactText \"716\u00d75=3580\" \"107\u00d77=749\"\nReplace-ExactText \"704\u00d79=6336\" \"520\u00d73=1560\"\nReplace-ExactText \"649\u00d74=2596\" \"795\u00d77=5565\"\nReplace-ExactText \"946\u00d77=6622\" \"192\u00d76=1152\"\nReplace-ExactText \"301\u00d79=2709\" \"850\u00d74=3400\"\nReplace-ExactText \"695\u00d72=1390\" \"200\u00d76=1200\"\nReplace-ExactText \"194\u00d75=970\" \"441\u00d75=2205\"\nReplace-ExactText \"506\u00d72=1012\" \"732\u00d79=6588\"\nReplace-ExactText \"588\u00d73=1764\" \"444\u00d74=1776\"\nReplace-ExactText \"140\u00d73=420\" \"542\u00d76=3252\"\nReplace-ExactText \"405\u00d79=3645\" \"773\u00d75=3865\"\nReplace-ExactText \"108\u00d77=756\" \"149\u00d73=447\"\nReplace-ExactText \"986\u00d72=1972\" \"167\u00d79=1503\"\nReplace-ExactText \"240\u00d76=1440\" \"171\u00d72=342\"\nReplace-ExactText \"562\u00d73=1686\" \"518\u00d75=2590\"\nReplace-ExactText \"481\u00d77=3367\" \"902\u00d72=1804\"\nReplace-ExactText \"716\u00d78=5728\" \"157\u00d76=942\"\n"}
